$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "B-javascript:S2201"
$ws.Range("A3").Value = "B-javascript:S2259"
$ws.Range("A4").Value = "B-javascript:S2583"
$ws.Range("A5").Value = "B-php:S836"
$ws.Range("A6").Value = "B-javascript:S3403"
$ws.Range("A7").Value = "B-php:S2201"
$ws.Range("A8").Value = "B-php:S1848"
$ws.Range("A9").Value = "B-php:S1656"
$ws.Range("A10").Value = "B-php:S3923"
$ws.Range("A11").Value = "V-php:S2964"
$ws.Range("A12").Value = "B-Web:BoldAndItalicTagsCheck"
$ws.Range("A13").Value = "B-javascript:S2757"
$ws.Range("A14").Value = "V-javascript:S2819"
$ws.Range("A15").Value = "B-php:S1763"
$ws.Range("A16").Value = "B-javascript:S2873"
$ws.Range("A17").Value = "B-javascript:S1143"
$ws.Range("A18").Value = "B-javascript:UnreachableCode"
$ws.Range("A19").Value = "V-php:S2053"
$ws.Range("A20").Value = "V-php:S2068"
$ws.Range("A21").Value = "B-javascript:S3981"
$ws.Range("A22").Value = "B-javascript:S3785"
$ws.Range("A23").Value = "B-php:S1764"
$ws.Range("A24").Value = "B-php:S1862"
$ws.Range("A25").Value = "B-javascript:S2137"
$ws.Range("A26").Value = "B-javascript:DuplicatePropertyName"
$ws.Range("A27").Value = "B-Web:DoctypePresenceCheck"
$ws.Range("A28").Value = "B-Web:FieldsetWithoutLegendCheck"
$ws.Range("A29").Value = "V-php:S4433"
$ws.Range("A30").Value = "B-php:S2757"
$ws.Range("A31").Value = "B-Web:PageWithoutTitleCheck"
$ws.Range("A32").Value = "B-javascript:S1656"
$ws.Range("A33").Value = "V-php:S4830"
$ws.Range("A34").Value = "B-Web:ImgWithoutAltCheck"
$ws.Range("A35").Value = "B-javascript:S905"
$ws.Range("A36").Value = "B-javascript:S3812"
$ws.Range("A37").Value = "B-javascript:S4043"
$ws.Range("A38").Value = "B-php:S1145"
$ws.Range("A39").Value = "B-php:S905"
$ws.Range("A40").Value = "B-Web:UnsupportedTagsInHtml5Check"
$ws.Range("A41").Value = "B-Web:FrameWithoutTitleCheck"
$ws.Range("A42").Value = "V-php:S4423"
